# Fix a mis-entered species code in the "assoc3" column (column R):
# the Google Drive update corrects "DIGL" -> "DRGL" for every data row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 3; $row -le 52; $row++) {
    $cell = $ws.Cells.Item($row, 18)   # column R = assoc3
    if ($cell.Value2 -eq "DIGL") {
        $cell.Value2 = "DRGL"
    }
}
